# Update "CPPbI" worksheet: split the existing single potential-share column
# into two columns - "energy related emissions" and "process emissions" -
# both holding the same values as the former column B, plus header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbI")

# New headers in row 1 (column A header "Unit: dimentionless ..." stays put)
$ws.Range("B1").Value = "energy related emissions"
$ws.Range("C1").Value = "process emissions"

# Mirror column B's existing values into the new column C, row by row.
for ($r = 2; $r -le 9; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value = $bVal
}

# Column widths as captured in the workbook after the edit (values chosen so
# the pixel-quantized stored width lands on the same bucket as the target).
$ws.Columns.Item(1).ColumnWidth = 42.834166666666746
$ws.Columns.Item(2).ColumnWidth = 23.834166666666633
$ws.Columns.Item(3).ColumnWidth = 25.0033333333333
